$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.271.48"
$ws.Range("E2").Value = "  +0.00%  "

# Row 3
$ws.Range("D3").Value = "3.916.62"
$ws.Range("E3").Value = "  -0.54%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.24%  "

# Row 5
$ws.Range("D5").Value = "'485.19"
$ws.Range("E5").Value = "  +1.42%  "

# Row 6
$ws.Range("D6").Value = "'146.94"
$ws.Range("E6").Value = "  -0.48%  "

# Row 7
$ws.Range("D7").Value = "'0.622"
$ws.Range("E7").Value = "  +0.24%  "

# Row 8
$ws.Range("D8").Value = "'0.998"
$ws.Range("E8").Value = "  -0.12%  "

# Row 9
$ws.Range("D9").Value = "'0.737"
$ws.Range("E9").Value = "  +1.05%  "

# Row 10
$ws.Range("E10").Value = "  +1.82%  "

# Row 11
$ws.Range("D11").Value = "'0.0000347"
$ws.Range("E11").Value = "  -1.36%  "

# Row 12
$ws.Range("D12").Value = "'43.18"
$ws.Range("E12").Value = "  +0.58%  "

# Row 13
$ws.Range("D13").Value = "'10.76"
$ws.Range("E13").Value = "  +3.70%  "

# Row 14
$ws.Range("D14").Value = "4.544.80"
$ws.Range("E14").Value = "  -0.79%  "

# Row 15
$ws.Range("D15").Value = "3.901.83"
$ws.Range("E15").Value = "  -2.84%  "

# Row 16
$ws.Range("D16").Value = "'14.27"
$ws.Range("E16").Value = "  -3.31%  "

# Row 17
$ws.Range("E17").Value = "  -0.50%  "

# Row 18
$ws.Range("E18").Value = "  +1.52%  "

# Row 19
$ws.Range("E19").Value = "  +0.86%  "

# Row 20
$ws.Range("D20").Value = "68.321.84"
$ws.Range("E20").Value = "  -0.04%  "

# Row 21
$ws.Range("D21").Value = "'431.12"
$ws.Range("E21").Value = "  -1.62%  "

# Row 22
$ws.Range("D22").Value = "'3.55"
$ws.Range("E22").Value = "  +6.98%  "

# Row 23
$ws.Range("D23").Value = "'15.06"
$ws.Range("E23").Value = "  +4.30%  "

# Row 24
$ws.Range("D24").Value = "'89.90"
$ws.Range("E24").Value = "  +2.74%  "

# Row 25
$ws.Range("D25").Value = "'11.62"
$ws.Range("E25").Value = "  +12.08%  "

# Row 26
$ws.Range("E26").Value = "  +2.36%  "

# Row 27
$ws.Range("D27").Value = "'11.23"
$ws.Range("E27").Value = "  +10.85%  "

# Row 28
$ws.Range("D28").Value = "'37.76"
$ws.Range("E28").Value = "  -1.31%  "

# Row 29
$ws.Range("E29").Value = "  -1.60%  "

# Row 30
$ws.Range("D30").Value = "'717.90"
$ws.Range("E30").Value = "  -1.36%  "

# Row 31
$ws.Range("D31").Value = "'13.76"
$ws.Range("E31").Value = "  +3.65%  "

# Row 32
$ws.Range("E32").Value = "  +2.12%  "

# Row 33
$ws.Range("E33").Value = "  +4.53%  "

# Row 34
$ws.Range("D34").Value = "0.0₃0890"
$ws.Range("E34").Value = "  +2.23%  "

# Row 35
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").Value = "'41.88"
$ws.Range("E35").Value = "  -0.85%  "

# Row 36
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "'6.14"
$ws.Range("E36").Value = "  +13.97%  "

# Row 37
$ws.Range("D37").Value = "'60.75"
$ws.Range("E37").Value = "  +2.08%  "

# Row 38
$ws.Range("D38").Value = "'0.402"
$ws.Range("E38").Value = "  +17.39%  "

# Row 39
$ws.Range("D39").Value = "'3.03"
$ws.Range("E39").Value = "  +6.88%  "

# Row 40
$ws.Range("E40").Value = "  -0.13%  "

# Row 41
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.144"
$ws.Range("E41").Value = "  -4.41%  "

# Row 42
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.0494"
$ws.Range("E42").Value = "  +4.95%  "

# Row 43
$ws.Range("E43").Value = "  +3.25%  "

# Row 44
$ws.Range("D44").Value = "'2.98"
$ws.Range("E44").Value = "  +2.67%  "

# Row 45
$ws.Range("E45").Value = "  +1.05%  "

# Row 46
$ws.Range("E46").Value = "  +5.36%  "

# Row 47
$ws.Range("E47").Value = "  +0.26%  "

# Row 48
$ws.Range("D48").Value = "'3.43"
$ws.Range("E48").Value = "  +0.81%  "

# Row 49
$ws.Range("D49").Value = "'2.14"
$ws.Range("E49").Value = "  -1.24%  "

# Row 50
$ws.Range("D50").Value = "'145.03"
$ws.Range("E50").Value = "  -1.36%  "

# Row 51
$ws.Range("E51").Value = "  +30.44%  "
